$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add Operacion / Zona columns, copying N1 header formatting ---
$ws.Range("O1").Value = "Operacion"
$ws.Range("P1").Value = "Zona"
$ws.Range("N1").Copy()
$ws.Range("O1:P1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Rows 2-54: Operacion / Zona data ---
$data1 = New-Object "object[,]" 53,2
$data1[0,0] = "Recoleta"
$data1[0,1] = "Capital Sur"
$data1[1,0] = "Recoleta"
$data1[1,1] = "Capital Sur"
$data1[2,0] = "Paternal"
$data1[2,1] = "Capital Norte"
$data1[3,0] = "Almagro"
$data1[3,1] = "Capital Sur"
$data1[4,0] = "Palermo"
$data1[4,1] = "Capital Sur"
$data1[5,0] = "Palermo"
$data1[5,1] = "Capital Sur"
$data1[6,0] = "San Telmo"
$data1[6,1] = "Capital Sur"
$data1[7,0] = "Palermo"
$data1[7,1] = "Capital Sur"
$data1[8,0] = "Almagro"
$data1[8,1] = "Capital Sur"
$data1[9,0] = "Almagro"
$data1[9,1] = "Capital Sur"
$data1[10,0] = "Boedo"
$data1[10,1] = "Capital Sur"
$data1[11,0] = "Colegiales"
$data1[11,1] = "Capital Norte"
$data1[12,0] = "Palermo"
$data1[12,1] = "Capital Sur"
$data1[13,0] = "Palermo"
$data1[13,1] = "Capital Sur"
$data1[14,0] = "Boedo"
$data1[14,1] = "Capital Sur"
$data1[15,0] = "Boedo"
$data1[15,1] = "Capital Sur"
$data1[16,0] = "San Telmo"
$data1[16,1] = "Capital Sur"
$data1[17,0] = "San Telmo"
$data1[17,1] = "Capital Sur"
$data1[18,0] = "Almagro"
$data1[18,1] = "Capital Sur"
$data1[19,0] = "Devoto"
$data1[19,1] = "Capital Norte"
$data1[20,0] = "Devoto"
$data1[20,1] = "Capital Norte"
$data1[21,0] = "San Telmo"
$data1[21,1] = "Capital Sur"
$data1[22,0] = "Recoleta"
$data1[22,1] = "Capital Sur"
$data1[23,0] = "Boedo"
$data1[23,1] = "Capital Sur"
$data1[24,0] = "Palermo"
$data1[24,1] = "Capital Sur"
$data1[25,0] = "Almagro"
$data1[25,1] = "Capital Sur"
$data1[26,0] = "Almagro"
$data1[26,1] = "Capital Sur"
$data1[27,0] = "Boedo"
$data1[27,1] = "Capital Sur"
$data1[28,0] = "San Telmo"
$data1[28,1] = "Capital Sur"
$data1[29,0] = "Colegiales"
$data1[29,1] = "Capital Norte"
$data1[30,0] = "Almagro"
$data1[30,1] = "Capital Sur"
$data1[31,0] = "Colegiales"
$data1[31,1] = "Capital Norte"
$data1[32,0] = "Boedo"
$data1[32,1] = "Capital Sur"
$data1[33,0] = "Almagro"
$data1[33,1] = "Capital Sur"
$data1[34,0] = "Boedo"
$data1[34,1] = "Capital Sur"
$data1[35,0] = "Almagro"
$data1[35,1] = "Capital Sur"
$data1[36,0] = "Palermo"
$data1[36,1] = "Capital Sur"
$data1[37,0] = "Boedo"
$data1[37,1] = "Capital Sur"
$data1[38,0] = "Colegiales"
$data1[38,1] = "Capital Norte"
$data1[39,0] = "Almagro"
$data1[39,1] = "Capital Sur"
$data1[40,0] = "Boedo"
$data1[40,1] = "Capital Sur"
$data1[41,0] = "Palermo"
$data1[41,1] = "Capital Sur"
$data1[42,0] = "Paternal"
$data1[42,1] = "Capital Norte"
$data1[43,0] = "Colegiales"
$data1[43,1] = "Capital Norte"
$data1[44,0] = "Colegiales"
$data1[44,1] = "Capital Norte"
$data1[45,0] = "Paternal"
$data1[45,1] = "Capital Norte"
$data1[46,0] = "Boedo"
$data1[46,1] = "Capital Sur"
$data1[47,0] = "Palermo"
$data1[47,1] = "Capital Sur"
$data1[48,0] = "Boedo"
$data1[48,1] = "Capital Sur"
$data1[49,0] = "Boedo"
$data1[49,1] = "Capital Sur"
$data1[50,0] = "Boedo"
$data1[50,1] = "Capital Sur"
$data1[51,0] = "Devoto"
$data1[51,1] = "Capital Norte"
$data1[52,0] = "Almagro"
$data1[52,1] = "Capital Sur"
$ws.Range("O2:P54").Value = $data1

# --- Rows 56-72: Operacion / Zona data ---
$data2 = New-Object "object[,]" 17,2
$data2[0,0] = "Almagro"
$data2[0,1] = "Capital Sur"
$data2[1,0] = "Almagro"
$data2[1,1] = "Capital Sur"
$data2[2,0] = "Almagro"
$data2[2,1] = "Capital Sur"
$data2[3,0] = "Almagro"
$data2[3,1] = "Capital Sur"
$data2[4,0] = "Almagro"
$data2[4,1] = "Capital Sur"
$data2[5,0] = "Almagro"
$data2[5,1] = "Capital Sur"
$data2[6,0] = "Devoto"
$data2[6,1] = "Capital Norte"
$data2[7,0] = "Palermo"
$data2[7,1] = "Capital Sur"
$data2[8,0] = "Colegiales"
$data2[8,1] = "Capital Norte"
$data2[9,0] = "Almagro"
$data2[9,1] = "Capital Sur"
$data2[10,0] = "Palermo"
$data2[10,1] = "Capital Sur"
$data2[11,0] = "Devoto"
$data2[11,1] = "Capital Norte"
$data2[12,0] = "Palermo"
$data2[12,1] = "Capital Sur"
$data2[13,0] = "Almagro"
$data2[13,1] = "Capital Sur"
$data2[14,0] = "Boedo"
$data2[14,1] = "Capital Sur"
$data2[15,0] = "Boedo"
$data2[15,1] = "Capital Sur"
$data2[16,0] = "Boedo"
$data2[16,1] = "Capital Sur"
$ws.Range("O56:P72").Value = $data2

